# Update ICAP_Nameplate (column F) values on Sheet1 for several generation
# types, as captured by the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 0.42    # Solar_PV
$ws.Range("F3").Value = 0.42    # Solar_Thermal
$ws.Range("F5").Value = 0.14699999999999999    # Wind_Onshore
$ws.Range("F12").Value = 0.42   # Res_Solar
$ws.Range("F13").Value = 0.42   # Comm_Solar
$ws.Range("F14").Value = 0.26   # Wind_Offshore_Shallow
$ws.Range("F15").Value = 0.26   # Wind_Offshore_Med
$ws.Range("F16").Value = 0.26   # Wind_Offshore_Deep

# Reflect the updated view state (frozen pane top-left cell and the active
# selection) recorded for the sheet after the edit.
$ws.Activate()
$ws.Range("D2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F17").Select()
